# Add the new row of data (Type3 / 7) that was appended to the bottom
# of the table on Sheet1, and leave the selection where Excel would
# naturally land after typing it in (the next empty cell, C16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Type3"
$ws.Range("B16").Value = 7

$ws.Range("C16").Select()
